# Actualización automática del inventario: agrega el nuevo producto
# "Rodillo de presion Hp" (código 7MD0M8) como la siguiente fila de la
# hoja de inventario (fila 5, justo debajo de los 3 productos existentes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 5
$ws.Range("A$row").Value = "7MD0M8"
$ws.Range("B$row").Value = "Rodillo de presion Hp"
$ws.Range("C$row").Value = "M107 M108 M2020"
$ws.Range("D$row").Value = 45000
$ws.Range("E$row").Value = 180000
$ws.Range("F$row").Value = 5
$ws.Range("G$row").Value = 8
$ws.Range("H$row").Value = "7MD0M8.jpg"
